$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Responsable" entries (shared string "Diego") for rows 6, 7, 9
$ws.Range("E6").Value = "Diego"
$ws.Range("E7").Value = "Diego"
$ws.Range("E9").Value = "Diego"

# Fill in the "Status" and "Responsable" for row 10, which was previously empty
$ws.Range("C10").Value = "ok"
$ws.Range("E10").Value = "Diego e Iván"

# Update the active selection to E10
$ws.Range("E10").Select()
